$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 84991.336
$ws.Range("J17").Value = 84991.336
$ws.Range("L17").Value = 254974.008
$ws.Range("N17").Value = -255310.008
$ws.Range("H18").Value = 2867.5
$ws.Range("I18").Value = 2490
$ws.Range("J18").Value = 4000
$ws.Range("K18").Value = 2490
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = -2206
$ws.Range("N18").Value = -4568
$ws.Range("H37").Value = 35132
$ws.Range("J37").Value = 250
$ws.Range("L37").Value = 750
$ws.Range("N37").Value = -1002
$ws.Range("H69").Value = 9582.333000000001
$ws.Range("I69").Value = 9123.25
$ws.Range("K69").Value = 27369.75
$ws.Range("M69").Value = -26495.75
$ws.Range("H72").Value = 9582.333000000001
$ws.Range("I72").Value = 9123.25
$ws.Range("K72").Value = 82109.25
$ws.Range("M72").Value = -77741.25
$ws.Range("H86").Value = 2873.2856
$ws.Range("I86").Value = 3033.6316
$ws.Range("J86").Value = 1350
$ws.Range("K86").Value = 3033.6316
$ws.Range("L86").Value = 1350
$ws.Range("M86").Value = -1910.6316
$ws.Range("N86").Value = -3596
$ws.Range("H89").Value = 2873.2856
$ws.Range("I89").Value = 3033.6316
$ws.Range("J89").Value = 1350
$ws.Range("K89").Value = 15168.158
$ws.Range("L89").Value = 6750
$ws.Range("M89").Value = -9552.158000000001
$ws.Range("N89").Value = -17982
$ws.Range("H106").Value = 8328.25
$ws.Range("I106").Value = 8328.25
$ws.Range("K106").Value = 8328.25
$ws.Range("M106").Value = -7697.25
$ws.Range("H137").Value = 2454.205
$ws.Range("I137").Value = 1744.6154
$ws.Range("J137").Value = 3873.3845
$ws.Range("K137").Value = 5233.8462
$ws.Range("L137").Value = 11620.1535
$ws.Range("M137").Value = -2683.8462
$ws.Range("N137").Value = -16720.1535
$ws.Range("H138").Value = 5429.3257
$ws.Range("I138").Value = 2715.2083
$ws.Range("J138").Value = 8857.684999999999
$ws.Range("K138").Value = 8145.624899999999
$ws.Range("L138").Value = 26573.055
$ws.Range("M138").Value = -3005.624899999999
$ws.Range("N138").Value = -36853.055

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9113.495999999999
$ws.Range("I32").Value = 9098.342000000001
$ws.Range("K32").Value = 9098.342000000001
$ws.Range("M32").Value = -8811.342000000001
$ws.Range("H74").Value = 1683.0385
$ws.Range("I74").Value = 1559.1305
$ws.Range("K74").Value = 1559.1305
$ws.Range("M74").Value = -685.1305
$ws.Range("H77").Value = 1683.0385
$ws.Range("I77").Value = 1559.1305
$ws.Range("K77").Value = 7795.6525
$ws.Range("M77").Value = -3427.6525
$ws.Range("H132").Value = 1726371
$ws.Range("I132").Value = 2128.9583
$ws.Range("K132").Value = 6386.874899999999
$ws.Range("M132").Value = -3856.874899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 39577.4
$ws.Range("I26").Value = 21944
$ws.Range("K26").Value = 21944
$ws.Range("M26").Value = -21652
$ws.Range("H80").Value = 1386.0344
$ws.Range("J80").Value = 1538.6666
$ws.Range("L80").Value = 1538.6666
$ws.Range("N80").Value = -3534.6666
$ws.Range("H83").Value = 1386.0344
$ws.Range("J83").Value = 1538.6666
$ws.Range("L83").Value = 7693.333000000001
$ws.Range("N83").Value = -17677.333
$ws.Range("H94").Value = 1493.6842
$ws.Range("I94").Value = 1706.3077
$ws.Range("J94").Value = 1033
$ws.Range("K94").Value = 1706.3077
$ws.Range("L94").Value = 1033
$ws.Range("M94").Value = -1255.3077
$ws.Range("N94").Value = -1935
$ws.Range("H105").Value = 590378.75
$ws.Range("I105").Value = 760431.4399999999
$ws.Range("J105").Value = 12199.6
$ws.Range("K105").Value = 760431.4399999999
$ws.Range("L105").Value = 12199.6
$ws.Range("M105").Value = -758684.4399999999
$ws.Range("N105").Value = -15693.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50004110
$ws.Range("I31").Value = 66669390
$ws.Range("J31").Value = 8282.200000000001
$ws.Range("K31").Value = 66669390
$ws.Range("L31").Value = 8282.200000000001
$ws.Range("M31").Value = -66669095
$ws.Range("N31").Value = -8872.200000000001
$ws.Range("H34").Value = 50004110
$ws.Range("I34").Value = 66669390
$ws.Range("J34").Value = 8282.200000000001
$ws.Range("K34").Value = 66669390
$ws.Range("L34").Value = 8282.200000000001
$ws.Range("M34").Value = -66669188
$ws.Range("N34").Value = -8686.200000000001
$ws.Range("H86").Value = 4843.0835
$ws.Range("I86").Value = 4931.091
$ws.Range("J86").Value = 3875
$ws.Range("K86").Value = 4931.091
$ws.Range("L86").Value = 3875
$ws.Range("M86").Value = -3808.091
$ws.Range("N86").Value = -6121
$ws.Range("H89").Value = 4843.0835
$ws.Range("I89").Value = 4931.091
$ws.Range("J89").Value = 3875
$ws.Range("K89").Value = 24655.455
$ws.Range("L89").Value = 19375
$ws.Range("M89").Value = -19039.455
$ws.Range("N89").Value = -30607
$ws.Range("H131").Value = 49326
$ws.Range("J131").Value = 49326
$ws.Range("L131").Value = 49326
$ws.Range("N131").Value = -59406
$ws.Range("H134").Value = 1532.081
$ws.Range("I134").Value = 1142.6
$ws.Range("K134").Value = 3427.8
$ws.Range("M134").Value = -892.7999999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 237.625
$ws.Range("I8").Value = 237.625
$ws.Range("K8").Value = 712.875
$ws.Range("M8").Value = -573.875
$ws.Range("H23").Value = 66666810
$ws.Range("J23").Value = 83333500
$ws.Range("L23").Value = 250000500
$ws.Range("N23").Value = -250000970
$ws.Range("H55").Value = 7575.3335
$ws.Range("I55").Value = 2750
$ws.Range("K55").Value = 8250
$ws.Range("M55").Value = -8073
$ws.Range("H113").Value = 1127.8462
$ws.Range("I113").Value = 1328.7142
$ws.Range("J113").Value = 1053.8422
$ws.Range("K113").Value = 3986.1426
$ws.Range("L113").Value = 3161.5266
$ws.Range("M113").Value = -1816.1426
$ws.Range("N113").Value = -7501.5266
$ws.Range("H117").Value = 3604.9285
$ws.Range("I117").Value = 213.7
$ws.Range("J117").Value = 12083
$ws.Range("K117").Value = 641.0999999999999
$ws.Range("L117").Value = 36249
$ws.Range("M117").Value = 2800.9
$ws.Range("N117").Value = -43133
$ws.Range("H131").Value = 4218.5483
$ws.Range("I131").Value = 3437.8333
$ws.Range("K131").Value = 10313.4999
$ws.Range("M131").Value = -5273.499899999999
$ws.Range("H139").Value = 4209.8857
$ws.Range("I139").Value = 1738.12
$ws.Range("J139").Value = 10389.3
$ws.Range("K139").Value = 5214.36
$ws.Range("L139").Value = 31167.9
$ws.Range("M139").Value = -74.35999999999967
$ws.Range("N139").Value = -41447.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 3175.6667
$ws.Range("J13").Value = 3689.9
$ws.Range("L13").Value = 3689.9
$ws.Range("N13").Value = -3967.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7084.3335
$ws.Range("I7").Value = 6702.8696
$ws.Range("J7").Value = 7759.231
$ws.Range("K7").Value = 6702.8696
$ws.Range("L7").Value = 7759.231
$ws.Range("M7").Value = -6590.8696
$ws.Range("N7").Value = -7983.231
$ws.Range("H82").Value = 3169.9285
$ws.Range("I82").Value = 1198.6666
$ws.Range("K82").Value = 1198.6666
$ws.Range("M82").Value = -837.6666
$ws.Range("H85").Value = 3169.9285
$ws.Range("I85").Value = 1198.6666
$ws.Range("K85").Value = 1198.6666
$ws.Range("M85").Value = 49.33339999999998
$ws.Range("H126").Value = 7084.3335
$ws.Range("I126").Value = 6702.8696
$ws.Range("J126").Value = 7759.231
$ws.Range("K126").Value = 20108.6088
$ws.Range("L126").Value = 23277.693
$ws.Range("M126").Value = -17638.6088
$ws.Range("N126").Value = -28217.693
$ws.Range("H140").Value = 241163.25
$ws.Range("J140").Value = 241163.25
$ws.Range("L140").Value = 241163.25
$ws.Range("N140").Value = -251523.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 53499.832
$ws.Range("J2").Value = 24500
$ws.Range("L2").Value = 24500
$ws.Range("N2").Value = -24724
$ws.Range("H113").Value = 646.09375
$ws.Range("I113").Value = 536.2917
$ws.Range("J113").Value = 975.5
$ws.Range("K113").Value = 1608.8751
$ws.Range("L113").Value = 2926.5
$ws.Range("M113").Value = 561.1249
$ws.Range("N113").Value = -7266.5
$ws.Range("H132").Value = 272947.38
$ws.Range("I132").Value = 2491.9
$ws.Range("J132").Value = 1432042.2
$ws.Range("K132").Value = 7475.700000000001
$ws.Range("L132").Value = 4296126.6
$ws.Range("M132").Value = -4945.700000000001
$ws.Range("N132").Value = -4301186.6
$ws.Range("H136").Value = 230355.45
$ws.Range("I136").Value = 3266.4324
$ws.Range("K136").Value = 9799.297200000001
$ws.Range("M136").Value = -7249.297200000001
